# Updates crypto price/volume data (and reorders a few coins) per the
# Aug 14 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.607.67'
$ws.Cells.Item(2, 5).Value = '  +0.85%  '

$ws.Cells.Item(3, 4).Value = '1.851.55'
$ws.Cells.Item(3, 5).Value = '  +0.18%  '

$ws.Cells.Item(4, 4).Value = '0.9988'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).Value = '240.63'
$ws.Cells.Item(5, 5).Value = '  +0.07%  '

$ws.Cells.Item(6, 4).Value = '0.6315'
$ws.Cells.Item(6, 5).Value = '  +0.50%  '

$ws.Cells.Item(7, 5).Value = '  -0.02%  '

$ws.Cells.Item(8, 5).Value = '  -1.38%  '

$ws.Cells.Item(9, 4).Value = '0.2919'
$ws.Cells.Item(9, 5).Value = '  +0.16%  '

$ws.Cells.Item(10, 4).Value = '25.14'
$ws.Cells.Item(10, 5).Value = '  +2.64%  '

$ws.Cells.Item(11, 4).Value = '0.07748'
$ws.Cells.Item(11, 5).Value = '  -0.06%  '

$ws.Cells.Item(12, 4).Value = '1.851.82'
$ws.Cells.Item(12, 5).Value = '  +0.21%  '

$ws.Cells.Item(13, 4).Value = '5.031'
$ws.Cells.Item(13, 5).Value = '  +0.47%  '

$ws.Cells.Item(14, 4).Value = '0.6845'
$ws.Cells.Item(14, 5).Value = '  +0.97%  '

$ws.Cells.Item(15, 4).Value = '0.00001028'
$ws.Cells.Item(15, 5).Value = '  -1.08%  '

$ws.Cells.Item(16, 4).Value = '82.79'
$ws.Cells.Item(16, 5).Value = '  -0.36%  '

$ws.Cells.Item(17, 4).Value = '6.345'
$ws.Cells.Item(17, 5).Value = '  +4.11%  '

$ws.Cells.Item(18, 4).Value = '29.590.53'
$ws.Cells.Item(18, 5).Value = '  +0.79%  '

$ws.Cells.Item(19, 4).Value = '230.74'
$ws.Cells.Item(19, 5).Value = '  +0.52%  '

$ws.Cells.Item(20, 4).Value = '12.41'
$ws.Cells.Item(20, 5).Value = '  +0.77%  '

$ws.Cells.Item(21, 4).Value = '1.000'
$ws.Cells.Item(21, 5).Value = '  +0.04%  '

$ws.Cells.Item(22, 4).Value = '7.543'
$ws.Cells.Item(22, 5).Value = '  +1.60%  '

$ws.Cells.Item(23, 5).Value = '  -0.11%  '

$ws.Cells.Item(24, 4).Value = '159.54'
$ws.Cells.Item(24, 5).Value = '  +0.24%  '

$ws.Cells.Item(25, 4).Value = '8.536'
$ws.Cells.Item(25, 5).Value = '  +1.13%  '

$ws.Cells.Item(26, 4).Value = '0.1369'
$ws.Cells.Item(26, 5).Value = '  -1.76%  '

$ws.Cells.Item(27, 4).Value = '17.58'
$ws.Cells.Item(27, 5).Value = '  -0.38%  '

$ws.Cells.Item(28, 2).Value = 'Toncoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(28, 4).Value = '1.477'
$ws.Cells.Item(28, 5).Value = '  +3.82%  '

$ws.Cells.Item(29, 2).Value = 'Hedera'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(29, 4).Value = '0.06600'
$ws.Cells.Item(29, 5).Value = '  +16.09%  '

$ws.Cells.Item(30, 4).Value = '1.488'
$ws.Cells.Item(30, 5).Value = '  +1.01%  '

$ws.Cells.Item(31, 4).Value = '4.115'
$ws.Cells.Item(31, 5).Value = '  -0.08%  '

$ws.Cells.Item(32, 4).Value = '4.110'
$ws.Cells.Item(32, 5).Value = '  +1.72%  '

$ws.Cells.Item(33, 5).Value = '  +1.66%  '

$ws.Cells.Item(34, 4).Value = '1.144'
$ws.Cells.Item(34, 5).Value = '  -0.97%  '

$ws.Cells.Item(35, 4).Value = '0.7003'
$ws.Cells.Item(35, 5).Value = '  +0.44%  '

$ws.Cells.Item(36, 4).Value = '2.568'
$ws.Cells.Item(36, 5).Value = '  -0.46%  '

$ws.Cells.Item(37, 4).Value = '0.01869'
$ws.Cells.Item(37, 5).Value = '  +2.23%  '

$ws.Cells.Item(38, 4).Value = '2.840'
$ws.Cells.Item(38, 5).Value = '  +4.59%  '

$ws.Cells.Item(39, 4).Value = '1.258.26'
$ws.Cells.Item(39, 5).Value = '  +1.62%  '

$ws.Cells.Item(40, 4).Value = '6.789'
$ws.Cells.Item(40, 5).Value = '  +5.83%  '

$ws.Cells.Item(41, 4).Value = '0.9355'
$ws.Cells.Item(41, 5).Value = '  +3.85%  '

$ws.Cells.Item(42, 5).Value = '  +0.19%  '

$ws.Cells.Item(43, 4).Value = '2.002.34'
$ws.Cells.Item(43, 5).Value = '  -0.09%  '

$ws.Cells.Item(44, 5).Value = '  +0.05%  '

$ws.Cells.Item(45, 4).Value = '66.37'
$ws.Cells.Item(45, 5).Value = '  +1.12%  '

$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(46, 4).Value = '7.124'
$ws.Cells.Item(46, 5).Value = '  -0.10%  '

$ws.Cells.Item(47, 2).Value = 'RenderToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(47, 4).Value = '1.736'
$ws.Cells.Item(47, 5).Value = '  +3.79%  '

$ws.Cells.Item(48, 4).Value = '0.1162'
$ws.Cells.Item(48, 5).Value = '  +0.72%  '

$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '9.022'
$ws.Cells.Item(49, 5).Value = '  +0.16%  '

$ws.Cells.Item(50, 2).Value = 'TheSandbox'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(50, 4).Value = '0.3949'
$ws.Cells.Item(50, 5).Value = '  -1.05%  '

$ws.Cells.Item(51, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(51, 4).Value = '0.00000000111'
$ws.Cells.Item(51, 5).Value = '  -2.80%  '
